$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 535.58185
$ws.Range("J17").Value = 535.58185
$ws.Range("L17").Value = 1606.74555
$ws.Range("N17").Value = -1942.74555

$ws.Range("H28").Value = 7843.385
$ws.Range("I28").Value = 878.1667
$ws.Range("K28").Value = 878.1667
$ws.Range("M28").Value = -393.1667

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H53").Value = 346.4
$ws.Range("J53").Value = 830.2
$ws.Range("L53").Value = 830.2
$ws.Range("N53").Value = -2104.2

$ws.Range("H61").Value = 3553
$ws.Range("I61").Value = 691.25
$ws.Range("K61").Value = 2073.75
$ws.Range("M61").Value = -1901.75

$ws.Range("H127").Value = 8917.6
$ws.Range("I127").Value = 8751.909
$ws.Range("J127").Value = 9373.25
$ws.Range("K127").Value = 26255.727
$ws.Range("L127").Value = 28119.75
$ws.Range("M127").Value = -21295.727
$ws.Range("N127").Value = -38039.75

$ws.Range("H137").Value = 9186.392
$ws.Range("I137").Value = 12619.2
$ws.Range("J137").Value = 8232.833000000001
$ws.Range("K137").Value = 37857.60000000001
$ws.Range("L137").Value = 24698.499
$ws.Range("M137").Value = -35307.60000000001
$ws.Range("N137").Value = -29798.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 19
$ws.Range("I38").Value = 19
$ws.Range("K38").Value = 19
$ws.Range("M38").Value = 448

$ws.Range("H45").Value = 3504.0833
$ws.Range("I45").Value = 2331.125
$ws.Range("K45").Value = 2331.125
$ws.Range("M45").Value = -1954.125

$ws.Range("H63").Value = 2301.3125
$ws.Range("I63").Value = 1420.1818
$ws.Range("J63").Value = 4239.8
$ws.Range("K63").Value = 1420.1818
$ws.Range("L63").Value = 4239.8
$ws.Range("M63").Value = -734.1818000000001
$ws.Range("N63").Value = -5611.8

$ws.Range("H66").Value = 2301.3125
$ws.Range("I66").Value = 1420.1818
$ws.Range("J66").Value = 4239.8
$ws.Range("K66").Value = 7100.909000000001
$ws.Range("L66").Value = 21199
$ws.Range("M66").Value = -3668.909000000001
$ws.Range("N66").Value = -28063

$ws.Range("H122").Value = 6988.615
$ws.Range("J122").Value = 7122
$ws.Range("L122").Value = 21366
$ws.Range("N122").Value = -26266

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8799.143
$ws.Range("I99").Value = 3699.5
$ws.Range("K99").Value = 3699.5
$ws.Range("M99").Value = -2201.5

$ws.Range("H105").Value = 37040724
$ws.Range("I105").Value = 142863360
$ws.Range("K105").Value = 142863360
$ws.Range("M105").Value = -142861613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3087.75
$ws.Range("I31").Value = 1499.2667
$ws.Range("K31").Value = 1499.2667
$ws.Range("M31").Value = -1204.2667

$ws.Range("H34").Value = 3087.75
$ws.Range("I34").Value = 1499.2667
$ws.Range("K34").Value = 1499.2667
$ws.Range("M34").Value = -1297.2667

$ws.Range("H62").Value = 14447.2
$ws.Range("I62").Value = 15211.714
$ws.Range("K62").Value = 15211.714
$ws.Range("M62").Value = -14587.714

$ws.Range("H65").Value = 14447.2
$ws.Range("I65").Value = 15211.714
$ws.Range("K65").Value = 76058.57000000001
$ws.Range("M65").Value = -72938.57000000001

$ws.Range("H99").Value = 6259.7144
$ws.Range("I99").Value = 3122
$ws.Range("K99").Value = 3122
$ws.Range("M99").Value = -1624

$ws.Range("H122").Value = 3393.84
$ws.Range("I122").Value = 2765.8823
$ws.Range("J122").Value = 4728.25
$ws.Range("K122").Value = 8297.6469
$ws.Range("L122").Value = 14184.75
$ws.Range("M122").Value = -5847.6469
$ws.Range("N122").Value = -19084.75

$ws.Range("H126").Value = 6259.7144
$ws.Range("I126").Value = 3122
$ws.Range("K126").Value = 9366
$ws.Range("M126").Value = -6896

$ws.Range("H134").Value = 16498.436
$ws.Range("I134").Value = 8448.637000000001
$ws.Range("K134").Value = 25345.911
$ws.Range("M134").Value = -22810.911

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2991.75
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2991.75
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 8975.25
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -9507.25

$ws.Range("H68").Value = 2000679.8
$ws.Range("I68").Value = 999.5
$ws.Range("J68").Value = 3333800
$ws.Range("K68").Value = 2998.5
$ws.Range("L68").Value = 10001400
$ws.Range("M68").Value = -2187.5
$ws.Range("N68").Value = -10003022

$ws.Range("H71").Value = 2000679.8
$ws.Range("I71").Value = 999.5
$ws.Range("J71").Value = 3333800
$ws.Range("K71").Value = 8995.5
$ws.Range("L71").Value = 30004200
$ws.Range("M71").Value = -4939.5
$ws.Range("N71").Value = -30012312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16063.091
$ws.Range("I80").Value = 3035
$ws.Range("J80").Value = 23507.715
$ws.Range("K80").Value = 3035
$ws.Range("L80").Value = 23507.715
$ws.Range("M80").Value = -2037
$ws.Range("N80").Value = -25503.715

$ws.Range("H83").Value = 16063.091
$ws.Range("I83").Value = 3035
$ws.Range("J83").Value = 23507.715
$ws.Range("K83").Value = 15175
$ws.Range("L83").Value = 117538.575
$ws.Range("M83").Value = -10183
$ws.Range("N83").Value = -127522.575

$ws.Range("H126").Value = 15430.286
$ws.Range("J126").Value = 15430.286
$ws.Range("L126").Value = 46290.858
$ws.Range("N126").Value = -51230.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4774.225
$ws.Range("I7").Value = 4875.5
$ws.Range("J7").Value = 4586.143
$ws.Range("K7").Value = 4875.5
$ws.Range("L7").Value = 4586.143
$ws.Range("M7").Value = -4763.5
$ws.Range("N7").Value = -4810.143

$ws.Range("H46").Value = 1173.3636
$ws.Range("I46").Value = 475
$ws.Range("J46").Value = 1572.4286
$ws.Range("K46").Value = 475
$ws.Range("L46").Value = 1572.4286
$ws.Range("M46").Value = -287
$ws.Range("N46").Value = -1948.4286

$ws.Range("H82").Value = 2993.3572
$ws.Range("J82").Value = 4634.8335
$ws.Range("L82").Value = 4634.8335
$ws.Range("N82").Value = -5356.8335

$ws.Range("H85").Value = 2993.3572
$ws.Range("J85").Value = 4634.8335
$ws.Range("L85").Value = 4634.8335
$ws.Range("N85").Value = -7130.8335

$ws.Range("H126").Value = 4774.225
$ws.Range("I126").Value = 4875.5
$ws.Range("J126").Value = 4586.143
$ws.Range("K126").Value = 14626.5
$ws.Range("L126").Value = 13758.429
$ws.Range("M126").Value = -12156.5
$ws.Range("N126").Value = -18698.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10017100
$ws.Range("I32").Value = 10017100
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 10017100
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -10016783
$ws.Range("N32").ClearContents()

$ws.Range("H62").Value = 8333.700000000001
$ws.Range("I62").Value = 4369.7144
$ws.Range("J62").Value = 17583
$ws.Range("K62").Value = 4369.7144
$ws.Range("L62").Value = 17583
$ws.Range("M62").Value = -3745.7144
$ws.Range("N62").Value = -18831

$ws.Range("H65").Value = 8333.700000000001
$ws.Range("I65").Value = 4369.7144
$ws.Range("J65").Value = 17583
$ws.Range("K65").Value = 21848.572
$ws.Range("L65").Value = 87915
$ws.Range("M65").Value = -18728.572
$ws.Range("N65").Value = -94155

$ws.Range("H122").Value = 3195
$ws.Range("I122").Value = 3195
$ws.Range("K122").Value = 9585
$ws.Range("M122").Value = -7135

$ws.Range("H132").Value = 4939.2114
$ws.Range("I132").Value = 2858.0967
$ws.Range("J132").Value = 19275.777
$ws.Range("K132").Value = 8574.2901
$ws.Range("L132").Value = 57827.33099999999
$ws.Range("M132").Value = -6044.2901
$ws.Range("N132").Value = -62887.33099999999
